# ---------------------------------------------------------------------------
# words_sources.xlsx update
#   - number every row across the 5 sheets sequentially in a new column A
#     (fruits already had this; vegetables/colors/animals/materials get it)
#   - materials sheet: add 3 new rows (hydrogen, nitrogen, oxygen), sorted
#     alphabetically with the rest, each with its own hyperlink
#   - materials!E3 (carbon dioxide image name) becomes a literal value
#     instead of a formula, because the file name now uses an underscore
#     instead of a space
#   - re-apply the (now bigger) materials range as the sheet's hidden
#     _FilterDatabase defined name (left behind by Data > Sort on the range)
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. vegetables (sheet2): add column A = 49..90
# ---------------------------------------------------------------------------
$wsVeg = $wb.Worksheets.Item("vegetables")
for ($r = 1; $r -le 42; $r++) {
    $wsVeg.Cells.Item($r, 1).Value2 = 48 + $r
}

# ---------------------------------------------------------------------------
# 2. colors (sheet3): add column A = 91..109
# ---------------------------------------------------------------------------
$wsCol = $wb.Worksheets.Item("colors")
for ($r = 1; $r -le 19; $r++) {
    $wsCol.Cells.Item($r, 1).Value2 = 90 + $r
}

# ---------------------------------------------------------------------------
# 3. animals (sheet4): add column A = 110..149
# ---------------------------------------------------------------------------
$wsAni = $wb.Worksheets.Item("animals")
for ($r = 1; $r -le 40; $r++) {
    $wsAni.Cells.Item($r, 1).Value2 = 109 + $r
}

# ---------------------------------------------------------------------------
# 4. materials (sheet5)
# ---------------------------------------------------------------------------
$wsMat = $wb.Worksheets.Item("materials")

# 4a. carbon dioxide row (currently row 3): filename now has an underscore,
#     so the old "=C3&".png"" formula no longer produces the right text --
#     replace it with the literal value.
$wsMat.Cells.Item(3, 5).Value2 = "carbon_dioxide.png"

# 4b. insert the 3 new rows, working from the bottom up so the row numbers
#     used below are never invalidated by an earlier insert.

# -- oxygen: inserted right before "paper" (currently row 14)
$wsMat.Rows.Item(14).Insert()
$wsMat.Cells.Item(14, 2).Value2 = 5
$wsMat.Cells.Item(14, 3).Value2 = "oxygen"
$wsMat.Cells.Item(14, 4).Value2 = "oxigén"
$wsMat.Cells.Item(14, 5).Formula = "=C14&"".png"""
$wsMat.Cells.Item(14, 6).Value2 = "https://www.pngegg.com/hu/png-wjgzd"
$wsMat.Hyperlinks.Add($wsMat.Cells.Item(14, 6), "https://www.pngegg.com/hu/png-wjgzd") | Out-Null
$wsMat.Cells.Item(14, 7).Value2 = 45026

# -- nitrogen: inserted right before the "oxygen" row we just created
$wsMat.Rows.Item(14).Insert()
$wsMat.Cells.Item(14, 2).Value2 = 5
$wsMat.Cells.Item(14, 3).Value2 = "nitrogen"
$wsMat.Cells.Item(14, 4).Value2 = "nitrogén"
$wsMat.Cells.Item(14, 5).Formula = "=C14&"".png"""
$wsMat.Cells.Item(14, 6).Value2 = "https://www.pngegg.com/hu/png-ygzrt"
$wsMat.Hyperlinks.Add($wsMat.Cells.Item(14, 6), "https://www.pngegg.com/hu/png-ygzrt") | Out-Null
$wsMat.Cells.Item(14, 7).Value2 = 45026

# -- hydrogen: inserted right before "ice" (currently row 11)
$wsMat.Rows.Item(11).Insert()
$wsMat.Cells.Item(11, 2).Value2 = 5
$wsMat.Cells.Item(11, 3).Value2 = "hydrogen"
$wsMat.Cells.Item(11, 4).Value2 = "hidrogén"
$wsMat.Cells.Item(11, 5).Formula = "=C11&"".png"""
$wsMat.Cells.Item(11, 6).Value2 = "https://www.pngegg.com/hu/png-mluvn"
$wsMat.Hyperlinks.Add($wsMat.Cells.Item(11, 6), "https://www.pngegg.com/hu/png-mluvn") | Out-Null
$wsMat.Cells.Item(11, 7).Value2 = 45026

# 4c. number column A for all 24 data rows (150..173)
for ($r = 1; $r -le 24; $r++) {
    $wsMat.Cells.Item($r, 1).Value2 = 149 + $r
}

# 4d. the trailing blank row (used to extend the range to row 25, matching
#     the sort/filter range below)
$wsMat.Cells.Item(25, 7).Value2 = $null

# 4e. re-create the hidden _FilterDatabase name over the new, larger range
#     (this is what Excel leaves behind after a Data > Sort on B1:G25)
$existingFilterNames = @($wsMat.Names)
foreach ($n in $existingFilterNames) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.Delete()
    }
}
$filterName = $wsMat.Names.Add("_xlnm._FilterDatabase", "=materials!`$B`$1:`$G`$25")
$filterName.Visible = $false

Write-Output "All edits applied"
